$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.536.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.737.44"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.83"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4922"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2668"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06298"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.732.61"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07056"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.72"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.604"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6119"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.49"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.0000"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007373"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +6.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.531.53"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.15%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.956.34"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.593"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.708"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.261"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.79"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.46"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "108.10"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.764"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.041"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08063"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.718"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04593"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.611"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.23%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.010"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.89%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6371"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8961"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.09%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.014"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.405"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.003"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01512"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.27"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.01%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.403"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.95%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3910"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.885"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1188"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05399"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.24%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.56"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.790"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.269"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.85"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.93%  "
